# Applies the "My second try" commit:
#  1. Inserts a new centered Title paragraph ("TEST002") at the very top of
#     the document, carrying the (relocated) "_GoBack" bookmark.
#  2. Cleans up the 8 "addr / wd" list items: drops the spell-check
#     <w:proofErr/> markers and collapses the split runs into the plain
#     run layout Word produces when the text is retyped as one string.
#  3. Stamps the FUNCTIONAL-simulation screenshot's <wp:inline> with the
#     wp14:anchorId / wp14:editId pair Word assigns when the drawing is
#     touched interactively.
#  4. Removes the stray "_GoBack" bookmark that used to sit alone in the
#     trailing empty paragraph (it now lives in the new title paragraph).

$d = $word.ActiveDocument

function Insert-RawXmlAt($range, $bodyInnerXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
           'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' +
           'xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" ' +
           'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' +
           'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" ' +
           'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' +
           '<w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# --- 1. New title paragraph "TEST002" with the _GoBack bookmark -----------
$titleXml = '<w:p><w:pPr><w:pStyle w:val="Title"/><w:jc w:val="center"/></w:pPr>' +
            '<w:r><w:t>TEST002</w:t></w:r>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Insert-RawXmlAt $d.Range(0, 0) $titleXml

# --- 2. Remove the now-orphaned _GoBack bookmark at the end of the doc ----
$d.Bookmarks.Item("_GoBack").Delete()

# --- 3. Collapse the 8 "addr ... wd = ..." list items into clean runs -----
$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr></w:pPr>'

function Set-AddrParagraph($paraIndex, $runsXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    Insert-RawXmlAt $p.Range ('<w:p>' + $listPPr + $runsXml + '</w:p>')
}

Set-AddrParagraph 7  ('<w:r><w:t xml:space="preserve">addr from 0 – 127 and wd = </w:t></w:r>' +
                       '<w:r><w:t>00000001000000100000001100000100</w:t></w:r>')
Set-AddrParagraph 8  ('<w:r><w:t xml:space="preserve">addr from 0 – 127 and wd = </w:t></w:r>' +
                       '<w:r><w:t>00000100000000110000001000000001</w:t></w:r>')
Set-AddrParagraph 9  ('<w:r><w:t xml:space="preserve">addr from 0 – 127 and wd = </w:t></w:r>' +
                       '<w:r><w:t>11111111111111101111110111111100</w:t></w:r>')
Set-AddrParagraph 10 ('<w:r><w:t>addr from 0 – 127 and wd =</w:t></w:r>' +
                       '<w:r><w:t xml:space="preserve"> 11111100111111011111111011111111</w:t></w:r>')
Set-AddrParagraph 11 ('<w:r><w:t>addr from 0 – 127 and wd = 00000001000000101111110111111100</w:t></w:r>')
Set-AddrParagraph 12 ('<w:r><w:t>addr from 0 – 127 and wd = 11111100111111010000001000000001</w:t></w:r>')
Set-AddrParagraph 13 ('<w:r><w:t>addr from 0 – 127 and wd = 00000100000000111111111011111111</w:t></w:r>')
Set-AddrParagraph 14 ('<w:r><w:t>addr from 0 – 127 and wd = 11111110111111110000001100000100</w:t></w:r>')

# --- 4. Stamp the screenshot drawing with wp14:anchorId / wp14:editId -----
$shp = $d.InlineShapes.Item(1)
$shpRange = $shp.Range
$shpRange.Text = ""   # clear the picture char but keep its paragraph intact

$drawingXml = '<w:p><w:r><w:rPr><w:noProof/></w:rPr><w:drawing>' +
  '<wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="4134F371" wp14:editId="276C14B0">' +
  '<wp:extent cx="6850380" cy="1150620"/><wp:effectExtent l="0" t="0" r="7620" b="0"/>' +
  '<wp:docPr id="2" name="Picture 2"/>' +
  '<wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr>' +
  '<a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
  '<pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/>' +
  '<pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr>' +
  '<pic:blipFill><a:blip r:embed="rId5" cstate="print">' +
  '<a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}">' +
  '<a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst>' +
  '</a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' +
  '<pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="6850380" cy="1150620"/></a:xfrm>' +
  '<a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr>' +
  '</pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'

$insPoint = $d.Paragraphs.Item(17).Range.Duplicate
$insPoint.Collapse(1)
Insert-RawXmlAt $insPoint $drawingXml

Write-Host "All edits applied."
